# ---------------------------------------------------------------------------
# Edit summary (per the canonical OOXML diff):
#   1. Slide 6's table changes its table style (tableStyleId) from the
#      custom "Table_0" style {8D5A9B55-5228-4842-874B-6D1E30AA7711} to the
#      built-in style {A8037AEB-15BB-4259-8966-5A0F3DC9050F}.
#   2. The presentation's theme colour scheme is switched from the
#      "Integral" palette to the default "Office Theme" palette (the two
#      embedded theme parts effectively swap their colour content).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{A8037AEB-15BB-4259-8966-5A0F3DC9050F}")

# --- 2. Theme colour scheme: Integral -> Office Theme ----------------------
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as COM RGB
# (0x00BBGGRR) integers of the target "Office Theme" palette.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}

# Best-effort rename (not persisted by every host, harmless if ignored).
try { $p.SlideMaster.Theme.Name = "Office Theme" } catch {}
